# version 2: Add Feature -> let manager to create schedule for department member
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column B width (matches the new "doc.doc.user_job_id" column text)
# NOTE: Excel quantizes ColumnWidth to whole pixels of the Normal-style font's
# max-digit-width (7px here), so 19.25 itself isn't reachable; 18.56 is the
# COM input that lands on the nearest attainable grid point (19.2857...).
$ws.Columns.Item(2).ColumnWidth = 18.56

# New mini-table describing the manager -> department member schedule feature
# (write order matches shared-string creation order: A4, C4, C5, B4)
$ws.Range("A4").Value = "Message_Manager"
$ws.Range("C4").Value = "AP_UNIT_ID"
$ws.Range("C5").Value = "A0X"
$ws.Range("B4").Value = "doc.doc.user_job_id"

$ws.Range("B7").HorizontalAlignment = -4131
$ws.Range("B7").VerticalAlignment = -4108

$ws.Range("B5").Value = 13
$ws.Range("B5").HorizontalAlignment = -4131
$ws.Range("B5").VerticalAlignment = -4160

# Fix wording of the "nvarchar(1) null or 1" note to quote the literal value
$ws.Range("Q2").Value = "nvarchar(1) null or '1'"

# Update view: scroll / selection moved
$ws.Application.ActiveWindow.ScrollColumn = 9
$ws.Range("Q2").Select()
